$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("D2").Value = "MAX_RUNTIME"
$ws1.Range("D3").Value = "lrs-200 JOB Failure"
$ws1.Range("D4").Value = "MAX_RUNTIME"
$ws1.Range("D3").Select()
$ws1.Activate()
